$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: TaskId 13->24, ContainerBoatID SG001->MY00202, BerthId 5->3,
# Status Completed->Scheduled, PublishTime updated, EndTime cleared.
$ws.Range("A2").Value = 24
$ws.Range("B2").Value = "MY00202"
$ws.Range("D2").Value = 3
$ws.Range("E2").Value = "Scheduled"
$ws.Range("F2").Value = 45369.99888888889
$ws.Range("H2").Value = ""
$ws.Range("H2").Style = "Normal"

# Row 3: TaskId 7->22, ContainerBoatID MY001->CB0003, BerthId 1->3,
# PublishTime cleared, StartTime populated, tug boat list trimmed.
$ws.Range("A3").Value = 22
$ws.Range("B3").Value = "CB0003"
$ws.Range("D3").Value = 3
$ws.Range("F3").Value = ""
$ws.Range("F3").Style = "Normal"
$ws.Range("G3").Value = 45369.99945601852
$ws.Range("G3").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("I3").Value = "NB002`n"
$ws.Rows.Item(3).AutoFit()
